$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1664298907891865
$ws.Range("C2").Value = -0.04869585134476913
$ws.Range("D2").Value = 0.5433751427531116
$ws.Range("B3").Value = -0.3597807503635829
$ws.Range("C3").Value = 0.5877681652748142
$ws.Range("B4").Value = 0.4493330185371963
$ws.Range("B5").Value = 0.3237572684721769
$ws.Range("C5").Value = 0.09851611570591239
$ws.Range("D5").Value = -0.1617696342879523
$ws.Range("E5").Value = -0.01038665294713449
$ws.Range("B6").Value = -0.05644442470388011
$ws.Range("C6").Value = 0.02554406905147947
$ws.Range("D6").Value = 0.01570249145987312
$ws.Range("B7").Value = -0.06258861287178845
$ws.Range("C7").Value = 0.009326612534822337
$ws.Range("B8").Value = -0.1500529593267638
$ws.Range("B9").Value = 0.09935271978382244
$ws.Range("C9").Value = -0.04610865354300479
$ws.Range("D9").Value = -0.1002609095359117
$ws.Range("E9").Value = -0.05204930396077505
$ws.Range("B10").Value = 0.001578648981956705
$ws.Range("C10").Value = -0.08555277690262475
$ws.Range("D10").Value = -0.07419436050613024
$ws.Range("B11").Value = -0.08576781653471441
$ws.Range("C11").Value = -0.07617959856547851
$ws.Range("B12").Value = -0.2246333953485248
$ws.Range("B13").Value = -0.7536610307181386
$ws.Range("C13").Value = 0.06605150356014075
$ws.Range("D13").Value = 0.04927257475788877
$ws.Range("B14").Value = -0.2232627506474126
$ws.Range("C14").Value = -0.1276150314848035
$ws.Range("B15").Value = 0.2871152093399901
